$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 1821
$ws.Range("K3").Value = 1731
$ws.Range("J4").Value = 432
$ws.Range("K4").Value = 380
$ws.Range("K5").Value = 116
$ws.Range("K6").Value = 2234
$ws.Range("J7").Value = 6240
$ws.Range("K7").Value = 6282
$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K4").Value = 27
$ws.Range("K6").Value = 51
$ws.Range("K7").Value = 175
$ws.Range("K8").Value = 412
$ws.Range("K10").Value = 38
$ws.Range("K11").Value = 134
$ws.Range("K15").Value = 58
$ws.Range("K16").Value = 14
$ws.Range("K19").Value = 171
$ws.Range("K20").Value = 138
$ws.Range("K22").Value = 16
$ws.Range("K23").Value = 57
$ws.Range("K27").Value = 69
$ws.Range("K29").Value = 303
$ws.Range("K31").Value = 71
$ws.Range("K33").Value = 255
$ws.Range("K35").Value = 9
$ws.Range("K36").Value = 72
$ws.Range("K37").Value = 209
$ws.Range("K40").Value = 13
$ws.Range("K42").Value = 217
$ws.Range("K46").Value = 12
$ws.Range("K53").Value = 99
$ws.Range("K54").Value = 102
$ws.Range("K57").Value = 15
$ws.Range("J63").Value = 35
$ws.Range("K63").Value = 20
$ws.Range("K66").Value = 28
$ws.Range("K67").Value = 243
$ws.Range("K69").Value = 17
$ws.Range("K71").Value = 17
$ws.Range("K73").Value = 63
$ws.Range("K76").Value = 91
$ws.Range("K77").Value = 42
$ws.Range("K78").Value = 81
$ws.Range("K79").Value = 168
$ws.Range("K83").Value = 131
$ws.Range("K85").Value = 313
$ws.Range("K86").Value = 43
$ws.Range("K88").Value = 81
$ws.Range("K89").Value = 81
$ws.Range("K90").Value = 55
$ws.Range("K92").Value = 29
$ws.Range("K94").Value = 74
$ws.Range("K97").Value = 54
$ws.Range("K99").Value = 119
$ws.Range("J101").Value = 6240
$ws.Range("K101").Value = 6282
$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K6").Value = 45
$ws.Range("K7").Value = 175
$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("K3").Value = 35
$ws.Range("K7").Value = 134
$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K4").Value = 12
$ws.Range("K6").Value = 27
$ws.Range("K7").Value = 81
$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K2").Value = 113
$ws.Range("K3").Value = 103
$ws.Range("K7").Value = 313
$ws = $wb.Worksheets.Item("Norwood Park")
$ws.Range("K2").Value = 4
$ws.Range("K7").Value = 17
$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("K6").Value = 54
$ws.Range("K7").Value = 99
$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 120
$ws.Range("K3").Value = 119
$ws.Range("K4").Value = 21
$ws.Range("K6").Value = 142
$ws.Range("K7").Value = 412
$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("K4").Value = 9
$ws.Range("K7").Value = 131
$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K2").Value = 73
$ws.Range("K3").Value = 96
$ws.Range("K7").Value = 255
$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K3").Value = 68
$ws.Range("K7").Value = 209
$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K2").Value = 39
$ws.Range("K7").Value = 119
$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("K6").Value = 30
$ws.Range("K7").Value = 71
$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K6").Value = 82
$ws.Range("K7").Value = 243
$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K3").Value = 34
$ws.Range("K6").Value = 40
$ws.Range("K7").Value = 102
$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K3").Value = 103
$ws.Range("K7").Value = 303
$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K5").Value = 10
$ws.Range("K7").Value = 171
$ws = $wb.Worksheets.Item("River North")
$ws.Range("K2").Value = 17
$ws.Range("K7").Value = 91
$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("K2").Value = 18
$ws.Range("K7").Value = 51
$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K2").Value = 53
$ws.Range("K6").Value = 97
$ws.Range("K7").Value = 217
$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("K2").Value = 12
$ws.Range("K7").Value = 38
$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("K2").Value = 26
$ws.Range("K7").Value = 81
$ws = $wb.Worksheets.Item("Jefferson Park")
$ws.Range("K2").Value = 5
$ws.Range("K7").Value = 12
$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("K6").Value = 14
$ws.Range("K7").Value = 57
$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K3").Value = 59
$ws.Range("K4").Value = 13
$ws.Range("K7").Value = 168
$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K3").Value = 42
$ws.Range("K7").Value = 138
$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("K2").Value = 29
$ws.Range("K3").Value = 26
$ws.Range("K7").Value = 72
$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("K6").Value = 34
$ws.Range("K7").Value = 74
$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("K3").Value = 13
$ws.Range("K7").Value = 58
$ws = $wb.Worksheets.Item("North Center")
$ws.Range("K5").Value = 14
$ws.Range("K6").Value = 28
$ws = $wb.Worksheets.Item("Gold Coast")
$ws.Range("K5").Value = 7
$ws.Range("K6").Value = 9
$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("K2").Value = 18
$ws.Range("K6").Value = 27
$ws.Range("K7").Value = 63
$ws = $wb.Worksheets.Item("West Town")
$ws.Range("K2").Value = 10
$ws.Range("K7").Value = 54
$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Range("K2").Value = 8
$ws.Range("K7").Value = 29
$ws = $wb.Worksheets.Item("United Center")
$ws.Range("K4").Value = 2
$ws.Range("K7").Value = 81
$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("K3").Value = 12
$ws.Range("K4").Value = 10
$ws.Range("K7").Value = 69
$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("K2").Value = 9
$ws.Range("K6").Value = 43
$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("K2").Value = 25
$ws.Range("K7").Value = 55
$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("K3").Value = 3
$ws.Range("K7").Value = 15
$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("K3").Value = 4
$ws.Range("K7").Value = 16
$ws = $wb.Worksheets.Item("Oakland")
$ws.Range("K6").Value = 4
$ws.Range("K7").Value = 17
$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("K2").Value = 21
$ws.Range("K7").Value = 42
$ws = $wb.Worksheets.Item("Hegewisch")
$ws.Range("K2").Value = 6
$ws.Range("K3").Value = 5
$ws.Range("K7").Value = 13
$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("K3").Value = 7
$ws.Range("K6").Value = 27
$ws = $wb.Worksheets.Item("Bucktown")
$ws.Range("K4").Value = 1
$ws.Range("K7").Value = 14
